$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new columns: "height" (F) and "weight" (G). Copy the header
# formatting from the existing E1 header cell so the new header cells
# match the bold/bordered/centered style used by the other headers.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:G1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("F1").Value2 = "height"
$ws.Range("G1").Value2 = "weight"

# The existing "fantasy points" values in column E move over to the new
# column G. Column E becomes the new "height" column (constant value) and
# column F becomes the new "weight" column (constant value).
for ($r = 2; $r -le 17; $r++) {
    $oldFantasyPoints = $ws.Cells.Item($r, 5).Value2

    $ws.Cells.Item($r, 5).Value2 = 6.333333333333333
    $ws.Cells.Item($r, 6).Value2 = 261
    $ws.Cells.Item($r, 7).Value2 = $oldFantasyPoints
}
